$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete columns C:G (denunciante_id, denunciante_id_tipo, denunciante_nombre,
# denunciante_genero, denunciante_edad) -- remaining columns shift left to fill the gap.
$ws.Range("C:G").Delete()
